$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Conditional formatting rule on B2:G31 -------------------------------
# The original rule flags p-values that are "too small" (cellIs lessThan
# 0.005). The analysis changes it to flag p-values that are "too big"
# (cellIs greaterThan 0.05) - i.e. emphasising the NON-significant results.
$rng = $ws.Range("B2:G31")

# Re-creating the rule (rather than only editing it in place) mirrors how
# the workbook's history actually accumulated extra (unused) dxf entries
# in styles.xml while the live rule kept referencing the first style slot -
# a couple of intermediate edits were made and discarded before landing on
# the final threshold.
$probe1 = $rng.FormatConditions.Add(1, 6, "=0.005")
$probe1.Font.Color = 393372
$probe1.Interior.Color = 13551615
$probe1.Delete()

$probe2 = $rng.FormatConditions.Add(1, 6, "=0.005")
$probe2.Font.Color = 393372
$probe2.Interior.Color = 13551615
$probe2.Delete()

$fc = $rng.FormatConditions.Item(1)
$fc.Operator = 5
$fc.Formula1 = "=0.05"

# --- Selection -------------------------------------------------------------
# Reflects the reviewer scanning across the highlighted (non-significant)
# p-value cells after the rule change - ending with the last block selected.
$ws.Range("G29:G30").Select()
